$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation with a leading apostrophe so numeric-looking
# strings (e.g. "1.00", "94.70") stay text instead of becoming Doubles,
# matching the original inlineStr/text cells in the workbook.

$ws.Range("D2").Value = "'44.275.51"
$ws.Range("E2").Value = "'  +0.64%  "
$ws.Range("D3").Value = "'2.241.75"
$ws.Range("E3").Value = "'  +0.01%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.03%  "
$ws.Range("D5").Value = "'307.37"
$ws.Range("E5").Value = "'  -2.59%  "
$ws.Range("D6").Value = "'94.70"
$ws.Range("E6").Value = "'  -4.44%  "
$ws.Range("E7").Value = "'  -0.42%  "
$ws.Range("D8").Value = "'1.01"
$ws.Range("E8").Value = "'  +0.30%  "
$ws.Range("E9").Value = "'  -1.33%  "
$ws.Range("D10").Value = "'34.66"
$ws.Range("E10").Value = "'  -4.29%  "
$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "'  -1.20%  "
$ws.Range("D12").Value = "'7.18"
$ws.Range("E12").Value = "'  -2.35%  "
$ws.Range("E13").Value = "'  +0.02%  "
$ws.Range("D14").Value = "'2.335.51"
$ws.Range("E14").Value = "'  +4.01%  "
$ws.Range("D15").Value = "'2.584.70"
$ws.Range("E15").Value = "'  +0.08%  "
$ws.Range("D16").Value = "'0.830"
$ws.Range("E16").Value = "'  -1.39%  "
$ws.Range("D17").Value = "'13.51"
$ws.Range("E17").Value = "'  -3.21%  "
$ws.Range("D18").Value = "'44.038.15"
$ws.Range("E18").Value = "'  +0.36%  "
$ws.Range("D19").Value = "'0.0₃0964"
$ws.Range("E19").Value = "'  -1.33%  "
$ws.Range("D20").Value = "'6.40"
$ws.Range("E20").Value = "'  +1.01%  "
$ws.Range("D21").Value = "'12.19"
$ws.Range("E21").Value = "'  -7.27%  "
$ws.Range("D22").Value = "'65.53"
$ws.Range("E22").Value = "'  -0.01%  "
$ws.Range("D23").Value = "'237.61"
$ws.Range("E23").Value = "'  +0.33%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "'  -0.97%  "
$ws.Range("D25").Value = "'1.99"
$ws.Range("E25").Value = "'  -1.37%  "
$ws.Range("D27").Value = "'38.79"
$ws.Range("E27").Value = "'  +6.60%  "
$ws.Range("B28").Value = "'Cosmos"
$ws.Range("C28").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'9.92"
$ws.Range("E28").Value = "'  -2.06%  "
$ws.Range("B29").Value = "'Toncoin"
$ws.Range("C29").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "'  +3.98%  "
$ws.Range("D30").Value = "'20.04"
$ws.Range("E30").Value = "'  -0.18%  "
$ws.Range("D31").Value = "'5.85"
$ws.Range("E31").Value = "'  -2.19%  "
$ws.Range("D32").Value = "'153.51"
$ws.Range("E32").Value = "'  -1.50%  "
$ws.Range("D33").Value = "'0.0794"
$ws.Range("E33").Value = "'  -5.30%  "
$ws.Range("E34").Value = "'  -1.74%  "
$ws.Range("D35").Value = "'3.15"
$ws.Range("E35").Value = "'  -3.81%  "
$ws.Range("E36").Value = "'  +1.81%  "
$ws.Range("E37").Value = "'  -0.19%  "
$ws.Range("E38").Value = "'  -7.02%  "
$ws.Range("E39").Value = "'  -0.63%  "
$ws.Range("D40").Value = "'3.81"
$ws.Range("E40").Value = "'  -4.40%  "
$ws.Range("D41").Value = "'14.31"
$ws.Range("E41").Value = "'  -7.17%  "
$ws.Range("E42").Value = "'  -2.68%  "
$ws.Range("E43").Value = "'  +0.13%  "
$ws.Range("D44").Value = "'1.747.94"
$ws.Range("E44").Value = "'  +2.92%  "
$ws.Range("D45").Value = "'82.95"
$ws.Range("E45").Value = "'  +0.21%  "
$ws.Range("E46").Value = "'  -1.97%  "
$ws.Range("B47").Value = "'Aave"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'99.62"
$ws.Range("E47").Value = "'  -2.02%  "
$ws.Range("B48").Value = "'THORChain"
$ws.Range("C48").Value = "'https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'4.93"
$ws.Range("E48").Value = "'  -4.62%  "
$ws.Range("B49").Value = "'Stacks"
$ws.Range("C49").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.59"
$ws.Range("E49").Value = "'  -0.74%  "
$ws.Range("D50").Value = "'8.08"
$ws.Range("E50").Value = "'  -0.58%  "
$ws.Range("D51").Value = "'54.74"
$ws.Range("E51").Value = "'  -2.81%  "
